$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old column C (NCTId), shifting
# NCTId..intervention_type (old C:L) right to D:M.
$ws.Range("C1").EntireColumn.Insert()

# New column header
$ws.Range("C1").Value2 = "statut_name"

# Map each row's status label (column B) to the French "statut_name" text
# and fill column C for the data rows (2-21).
$labels = @{
    "noir"   = "pas de résultat ni de publication"
    "vert"   = "résultat et / ou publication posté dans les 12 mois"
    "orange" = "résultat et / ou publication posté dans les 36 mois"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $label = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 3).Value2 = $labels[$label]
}
